# "Tests with exit() function" — bump the staked-token amount for the
# second 150-day-bucket row (G14) from 7E+19 to 7E+26 and leave the
# selection where the user last clicked (G15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KWENTA")
$ws.Activate()

$ws.Range("G14").Value = [double]"7E+26"

$ws.Range("G15").Select() | Out-Null
